$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append three more seasons below the existing data (rows 2-6 already hold
# 2021-2017); continue the year sequence down through 2014.
$ws.Range("A7").Value = 2016
$ws.Range("A8").Value = 2015
$ws.Range("A9").Value = 2014

# Move the active selection from C10 to A10.
[void]$ws.Range("A10").Select()
